$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 previously held a log entry (date string in A7, minutes in B7).
# Remove that entry: B7 is fully cleared (cell removed), A7 keeps its
# existing center-aligned style but loses its value.
$ws.Range("B7").Clear()
$ws.Range("A7").ClearContents()

# Move the active cell selection to C10
$ws.Range("C10").Select()
